# Auto-generated edit script applying the Bahamut_Profits.xlsx diff
# Updates currentAveragePrice / currentAveragePriceNQ/HQ / LevePriceNQ/HQ / LeveProfitNQ/HQ
# figures (columns H-N) for specific leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 848.2857
$ws.Range("I137").Value = 800.7778
$ws.Range("J137").Value = 898.58826
$ws.Range("K137").Value = 2402.3334
$ws.Range("L137").Value = 2695.76478
$ws.Range("M137").Value = 147.6666
$ws.Range("N137").Value = -7795.76478

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 500
$ws.Range("K14").Value = 500
$ws.Range("M14").Value = -325

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1471.3
$ws.Range("I94").Value = 1436.8235
$ws.Range("J94").Value = 1666.6666
$ws.Range("K94").Value = 1436.8235
$ws.Range("L94").Value = 1666.6666
$ws.Range("M94").Value = -985.8235
$ws.Range("N94").Value = -2568.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 770.3077
$ws.Range("I16").Value = 680.1429000000001
$ws.Range("J16").Value = 875.5
$ws.Range("K16").Value = 680.1429000000001
$ws.Range("L16").Value = 875.5
$ws.Range("M16").Value = -393.1429000000001
$ws.Range("N16").Value = -1449.5

$ws.Range("H31").Value = 2708
$ws.Range("I31").Value = 2848.7334
$ws.Range("J31").Value = 2485.7896
$ws.Range("K31").Value = 2848.7334
$ws.Range("L31").Value = 2485.7896
$ws.Range("M31").Value = -2553.7334
$ws.Range("N31").Value = -3075.7896

$ws.Range("H34").Value = 2708
$ws.Range("I34").Value = 2848.7334
$ws.Range("J34").Value = 2485.7896
$ws.Range("K34").Value = 2848.7334
$ws.Range("L34").Value = 2485.7896
$ws.Range("M34").Value = -2646.7334
$ws.Range("N34").Value = -2889.7896

$ws.Range("H58").Value = 1064.1794
$ws.Range("I58").Value = 1080.6451
$ws.Range("J58").Value = 1000.375
$ws.Range("K58").Value = 1080.6451
$ws.Range("L58").Value = 1000.375
$ws.Range("M58").Value = -877.6451
$ws.Range("N58").Value = -1406.375

$ws.Range("H94").Value = 2585.1052
$ws.Range("I94").Value = 3074.25
$ws.Range("J94").Value = 2454.6667
$ws.Range("K94").Value = 3074.25
$ws.Range("L94").Value = 2454.6667
$ws.Range("M94").Value = -2623.25
$ws.Range("N94").Value = -3356.6667

$ws.Range("H99").Value = 2309.279
$ws.Range("I99").Value = 1889.6552
$ws.Range("J99").Value = 3178.5
$ws.Range("K99").Value = 1889.6552
$ws.Range("L99").Value = 3178.5
$ws.Range("M99").Value = -391.6551999999999
$ws.Range("N99").Value = -6174.5

$ws.Range("H113").Value = 770.3077
$ws.Range("I113").Value = 680.1429000000001
$ws.Range("J113").Value = 875.5
$ws.Range("K113").Value = 680.1429000000001
$ws.Range("L113").Value = 875.5
$ws.Range("M113").Value = 1489.8571
$ws.Range("N113").Value = -5215.5

$ws.Range("H122").Value = 947
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 2309.279
$ws.Range("I126").Value = 1889.6552
$ws.Range("J126").Value = 3178.5
$ws.Range("K126").Value = 5668.9656
$ws.Range("L126").Value = 9535.5
$ws.Range("M126").Value = -3198.9656
$ws.Range("N126").Value = -14475.5

$ws.Range("H136").Value = 1064.1794
$ws.Range("I136").Value = 1080.6451
$ws.Range("J136").Value = 1000.375
$ws.Range("K136").Value = 3241.9353
$ws.Range("L136").Value = 3001.125
$ws.Range("M136").Value = -691.9353000000001
$ws.Range("N136").Value = -8101.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1755
$ws.Range("I25").Value = 10
$ws.Range("J25").Value = 3500
$ws.Range("K25").Value = 30
$ws.Range("L25").Value = 10500
$ws.Range("M25").Value = 139
$ws.Range("N25").Value = -10838

$ws.Range("H30").Value = 1755
$ws.Range("I30").Value = 10
$ws.Range("J30").Value = 3500
$ws.Range("K30").Value = 30
$ws.Range("L30").Value = 10500
$ws.Range("M30").Value = 72
$ws.Range("N30").Value = -10704

$ws.Range("H113").Value = 596.08826
$ws.Range("I113").Value = 793
$ws.Range("J113").Value = 514.0417
$ws.Range("K113").Value = 2379
$ws.Range("L113").Value = 1542.1251
$ws.Range("M113").Value = -209
$ws.Range("N113").Value = -5882.1251

$ws.Range("H137").Value = 30711.361
$ws.Range("J137").Value = 56556.21
$ws.Range("L137").Value = 169668.63
$ws.Range("N137").Value = -179868.63

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 67004.8
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 67004.8
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 67004.8
$ws.Range("N29").Value = -67584.8
$ws.Range("M29").ClearContents()

$ws.Range("H31").Value = 5744.2
$ws.Range("I31").Value = 930.25
$ws.Range("K31").Value = 930.25
$ws.Range("M31").Value = -638.25

$ws.Range("H37").Value = 5744.2
$ws.Range("I37").Value = 930.25
$ws.Range("K37").Value = 930.25
$ws.Range("M37").Value = -653.25

$ws.Range("H80").Value = 3584.647
$ws.Range("I80").Value = 3569.6667
$ws.Range("J80").Value = 3697
$ws.Range("K80").Value = 3569.6667
$ws.Range("L80").Value = 3697
$ws.Range("M80").Value = -2571.6667
$ws.Range("N80").Value = -5693

$ws.Range("H83").Value = 3584.647
$ws.Range("I83").Value = 3569.6667
$ws.Range("J83").Value = 3697
$ws.Range("K83").Value = 17848.3335
$ws.Range("L83").Value = 18485
$ws.Range("M83").Value = -12856.3335
$ws.Range("N83").Value = -28469

$ws.Range("H102").Value = 1249.909
$ws.Range("I102").Value = 1216.9
$ws.Range("J102").Value = 1580
$ws.Range("K102").Value = 1216.9
$ws.Range("L102").Value = 1580
$ws.Range("M102").Value = 405.0999999999999
$ws.Range("N102").Value = -4824

$ws.Range("H132").Value = 1751
$ws.Range("I132").Value = 1777.4546
$ws.Range("J132").Value = 1696.4375
$ws.Range("K132").Value = 5332.3638
$ws.Range("L132").Value = 5089.3125
$ws.Range("M132").Value = -2802.3638
$ws.Range("N132").Value = -10149.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14:N14").ClearContents()

$ws.Range("H122").Value = 6450.636
$ws.Range("I122").Value = 7671.4116
$ws.Range("K122").Value = 23014.2348
$ws.Range("M122").Value = -20564.2348

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 39000
$ws.Range("J70").Value = 39000
$ws.Range("L70").Value = 39000
$ws.Range("N70").Value = -39630

$ws.Range("H73").Value = 39000
$ws.Range("J73").Value = 39000
$ws.Range("L73").Value = 39000
$ws.Range("N73").Value = -41184
